$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ACTIVITY DIAGRAM section (row 46) is now complete: give its header
# row the same "done" highlight fill already used by the other completed
# section headers (USE CASE DIAGRAM / CLASS DIAGRAM, row 3 / row 13), and
# extend the highlight into column C the same way those rows do.
$ws.Range("A3:B3").Copy()
$ws.Range("A46:B46").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C46").PasteSpecial(-4122)

# --- STATE MACHINE DIAGRAM section still needs work - leave a note.
$ws.Range("C25").Value = "state machine diagram da sistemare"

# --- Mark newly completed checklist items ("x") in the ACTIVITY DIAGRAM
# section, and clear one that moved.
$ws.Range("C47").Value = "x"
$ws.Range("C50").Value = "x"
$ws.Range("C51").Value = "x"
$ws.Range("C52").ClearContents()
$ws.Range("C55").Value = "x"
$ws.Range("C56").Value = "x"

# --- Restore the selection / scroll position to where the user left off.
$ws.Range("C52").Select()
$excel.ActiveWindow.ScrollRow = 35
